# Updates the cryptos price/volume table (cols B-E) on Sheet1 to match the
# latest scrape. Cells in column D that look like plain numbers ("40.81",
# "0.5110", ...) are stored as TEXT in the workbook (t="inlineStr"), so a
# leading apostrophe is used to stop Excel's COM layer from auto-coercing
# them to numeric values; ".Style" is reset to "Normal" right after so no
# extra cell formatting (quote-prefix) is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.735.43"
$ws.Range("E2").Value = "  -2.75%  "
$ws.Range("D3").Value = "1.783.53"
$ws.Range("E3").Value = "  -2.29%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'310.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.13%  "
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").Value = "'0.5110"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.54%  "
$ws.Range("D8").Value = "'0.3870"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.60%  "
$ws.Range("D9").Value = "'0.07811"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -7.38%  "
$ws.Range("B10").Value = "OKB"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D10").Value = "'40.81"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.62%  "
$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").Value = "'1.086"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.59%  "
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("D13").Value = "'6.206"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.39%  "
$ws.Range("D14").Value = "'20.15"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.51%  "
$ws.Range("D15").Value = "1.779.90"
$ws.Range("E15").Value = "  -2.38%  "
$ws.Range("D16").Value = "'7.191"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.50%  "
$ws.Range("D17").Value = "'91.30"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.91%  "
$ws.Range("D18").Value = "'0.00001075"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.23%  "
$ws.Range("D19").Value = "'0.06549"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.79%  "
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").Value = "'16.99"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.45%  "
$ws.Range("D22").Value = "'5.898"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.84%  "
$ws.Range("D23").Value = "27.785.72"
$ws.Range("E23").Value = "  -2.65%  "
$ws.Range("D24").Value = "'10.97"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.06%  "
$ws.Range("D25").Value = "'2.224"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.31%  "
$ws.Range("D26").Value = "'160.42"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.57%  "
$ws.Range("B27").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C27").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D27").Value = "1.985.59"
$ws.Range("E27").Value = "  -2.25%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'20.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.59%  "
$ws.Range("D29").Value = "'2.359"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.77%  "
$ws.Range("D30").Value = "'123.38"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.76%  "
$ws.Range("D31").Value = "'0.1077"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.34%  "
$ws.Range("E32").Value = "  -6.39%  "
$ws.Range("D33").Value = "'3.633"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.74%  "
$ws.Range("D34").Value = "'5.470"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.50%  "
$ws.Range("D35").Value = "'0.07032"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.43%  "
$ws.Range("D36").Value = "'0.02302"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.62%  "
$ws.Range("D37").Value = "'8.744"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.87%  "
$ws.Range("E38").Value = "  -4.85%  "
$ws.Range("B39").Value = "Aptos"
$ws.Range("C39").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D39").Value = "'11.44"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.50%  "
$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").Value = "'4.976"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.86%  "
$ws.Range("D41").Value = "'0.6073"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.04%  "
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("D43").Value = "'1.146"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.93%  "
$ws.Range("B44").Value = "WEMIXTOKEN"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "'1.320"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.67%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'13.05"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.65%  "
$ws.Range("D46").Value = "'0.5875"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.66%  "
$ws.Range("D47").Value = "'3.696"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.36%  "
$ws.Range("D48").Value = "'124.43"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.76%  "
$ws.Range("D49").Value = "'1.200"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.54%  "
$ws.Range("D50").Value = "'1.897"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Value = "'0.06835"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.16%  "
